$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows to grow the table from 15 to 18 rows --------------
# A new "Condition" row is inserted ahead of each book row (except the very
# first one, "Introduction", which already had its own condition row built
# in), and a brand new book ("Table Processing") is introduced.
$ws.Rows("3:3").Insert()
$ws.Rows("6:6").Insert()
$ws.Rows("9:9").Insert()

# --- Row 3: new stand-alone "Inputs_Index" condition row ------------------
$ws.Range("D3").Clear()
$ws.Range("E3").Value = "Inputs_Index\*.txt"

# --- Row 6: "B4P Language" renamed to "Language Guide" --------------------
$ws.Range("B6").Value = "Language Guide"
$ws.Range("C6").Value = "LAN"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = "Inputs_Language\*.txt"

# --- Row 7: stand-alone "Inputs_Index" condition row (was the book row) ---
$ws.Range("B7").Clear()
$ws.Range("C7").Clear()
$ws.Range("D7").Clear()
$ws.Range("E7").Value = "Inputs_Index\*.txt"

# --- Row 8: brand-new book "Table Processing" / "TAB" ---------------------
$ws.Range("B8").Value = "Table Processing"
$ws.Range("C8").Value = "TAB"
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = "Inputs_Tables\*.txt"

# --- Row 9: stand-alone "Inputs_Index" condition row -----------------------
$ws.Range("E9").Value = "Inputs_Index\*.txt"

# --- Row 10: "Function Library Guide" renamed to "Function Library" -------
$ws.Range("B10").Value = "Function Library"
$ws.Range("D10").Value = 5

# --- Update the saved selection --------------------------------------------
$ws.Range("B9").Select()
